$d = $word.ActiveDocument

# 1. Cover page date: "08 April 2013" -> "12 June 2013"
$d.Content.Find.Execute("April", $true, $true, $false, $false, $false, $true, 1, $false, "June", 2) | Out-Null
$d.Content.Find.Execute("08", $true, $true, $false, $false, $false, $true, 1, $false, "12", 2) | Out-Null

# 2. "svn co http://epicsqt.svn.sourceforge.net/svnroot/epicsqt/trunk " -> "svn checkout svn://svn.code.sf.net/p/epicsqt/code/trunk "
$d.Content.Find.Execute(" co http://epicsqt.svn.sourceforge.net/svnroot/epicsqt/trunk ", $true, $false, $false, $false, $false, $true, 1, $false, " checkout svn://svn.code.sf.net/p/epicsqt/code/trunk ", 2) | Out-Null

# 3. viewvc browse link display text
$d.Content.Find.Execute("http://epicsqt.svn.sourceforge.net/viewvc/epicsqt/", $true, $false, $false, $false, $false, $true, 1, $false, "http://sourceforge.net/p/epicsqt/code/HEAD/tree/", 2) | Out-Null

# 4. "svn export ..." makefile command
$d.Content.Find.Execute(" export http://epicsqt.svn.sourceforge.net/svnroot/epicsqt/trunk/resources/makefile", $true, $false, $false, $false, $false, $true, 1, $false, " export http://svn.code.sf.net/p/epicsqt/code/trunk/resources/makefile", 2) | Out-Null
